$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells F1, G1, H1 - copy style from existing header (e.g. E1)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null

# New boolean data cells F2:H4 set to FALSE
$ws.Range("F2:H4").Value = $false
